$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "RegLog"

# --- Update RegLog's "Decisao" (D) column formula: new threshold (0.6) and numeric 1/0 output
#     instead of the text "Aprovado"/"Reprovado" (threshold 0.75) ---
for ($r = 2; $r -le 16; $r++) {
    $ws1.Cells.Item($r, 4).Formula = "=IF(C$r>=0.6,1,0)"
}

# Preserve RegLog's own selection (D18) before we switch sheets away from it -
# this engine keeps only the selection of whichever sheet was active at the
# time .Select() was called.
$ws1.Range("D18").Select()

# --- Duplicate RegLog (brings along column widths, styles, conditional
#     formatting, autofilter and the embedded chart) to seed the new "blank"
#     sheet, then strip it back down to a blank template ---
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "blank"

# Reset B:D to empty (keep styles/number formats, drop formulas+values)
$ws2.Range("B2:D16").ClearContents()

# Restore the original (unsorted) raw "Nota Simulado" sample for the blank
# template's A column
$notas = @(7,3,8,7,9,7,5,2,8,8,3,5,6,9,7)
for ($i = 0; $i -lt $notas.Length; $i++) {
    $ws2.Cells.Item($i + 2, 1).Value = $notas[$i]
}

# The sheet-copy also duplicated the chart as a *linked* clone (same
# underlying chart object as RegLog's) rather than an independent one, so
# drop it and build a fresh scatter chart sourced from the blank sheet.
if ($ws2.ChartObjects().Count -gt 0) {
    $ws2.ChartObjects().Item(1).Delete()
}

$co2 = $ws2.ChartObjects().Add(346275, 76201, 2247900, 2247900)
$chart2 = $co2.Chart
$chart2.ChartType = -4169
$ser2 = $chart2.SeriesCollection().NewSeries()
$ser2.Name = "=blank!`$C`$1"
$ser2.XValues = $ws2.Range("A2:A16")
$ser2.Values = $ws2.Range("C2:C16")
$chart2.HasTitle = $true
$chart2.ChartTitle.Text = "Função sigmoide"

# --- Activate "blank" as the visible tab with its own selection ---
$ws2.Select()
$ws2.Range("B2").Select()
